# fix: add domain field in EU funds report
#
# Adds a new "Domain" header column (M1) to the "Data" sheet, right after the
# existing "Postal Address" column (L1), extends the column widths/autofilter/
# filter-database range to cover the new column, and updates the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Add the new "Domain" header cell in M1 --------------------------------
# Copy the formatting from L1 ("Postal Address") so the new header cell gets
# the same header style (fill/etc.) as the rest of the header row.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M1").Value = "Domain"

# --- Column widths -----------------------------------------------------
# Previously: K (11) = 18.5, L (12) = 22.5
# Now:        K (11) = 18.5, L (12) = 18.5, M (13, new) = 22.5
# (COM ColumnWidth is offset from the raw stored XML width by 5/6, so we
# subtract that offset to land on the exact target widths.)
$ws.Columns.Item(12).ColumnWidth = 18.5 - (5/6)
$ws.Columns.Item(13).ColumnWidth = 22.5 - (5/6)

# --- Extend the AutoFilter range to include the new column -----------------
$ws.AutoFilterMode = $false
$ws.Range("A1:M1").AutoFilter()

# --- Keep the hidden _FilterDatabase defined name in sync ------------------
try {
    $fdb = $wb.Names.Item("Data!_FilterDatabase")
    $fdb.RefersTo = "=Data!`$A`$1:`$M`$1"
} catch {
    $wb.Names.Add("Data!_FilterDatabase", "=Data!`$A`$1:`$M`$1")
}

# --- Update the selected cell in the sheet ----------------------------------
$ws.Range("J10").Select()
